# Generate Report for Handoff
#
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the file 533db0e0-4a0d-4377-99c6-29b43b86da7e.md (row 6)
# on the Overview, zh-cn and de-de sheets, reflecting a fresh handoff
# report generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-16 00:38:30"

# --- zh-cn sheet: column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-16 00:38:25"

# --- de-de sheet: column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-16 00:38:30"
